$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '47.475.29'
$ws.Range("E2").Value = '  +4.65%  '
$ws.Range("D3").Value = '2.492.46'
$ws.Range("E3").Value = '  +2.82%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '322.83'
$ws.Range("E5").Value = '  +1.25%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '105.63'
$ws.Range("E6").Value = '  +2.41%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.525'
$ws.Range("E7").Value = '  +1.67%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.999'
$ws.Range("E8").Value = '  -0.03%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.545'
$ws.Range("E9").Value = '  +2.97%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '38.24'
$ws.Range("E10").Value = '  +7.48%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0813'
$ws.Range("E11").Value = '  +1.13%  '
$ws.Range("E12").Value = '  +1.16%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '18.41'
$ws.Range("E13").Value = '  +1.01%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.19'
$ws.Range("E14").Value = '  +1.63%  '
$ws.Range("D15").Value = '2.878.95'
$ws.Range("E15").Value = '  +2.73%  '
$ws.Range("D16").Value = '2.490.25'
$ws.Range("E16").Value = '  +3.15%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.849'
$ws.Range("E17").Value = '  +0.70%  '
$ws.Range("D18").Value = '47.361.28'
$ws.Range("E18").Value = '  +4.56%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.80'
$ws.Range("E19").Value = '  +4.79%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.58'
$ws.Range("E20").Value = '  +3.70%  '
$ws.Range("D21").Value = '0.0₃0938'
$ws.Range("E21").Value = '  +1.55%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '70.72'
$ws.Range("E22").Value = '  +2.57%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.42'
$ws.Range("E23").Value = '  +6.43%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '251.55'
$ws.Range("E24").Value = '  +2.82%  '
$ws.Range("E25").Value = '  +3.27%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '26.27'
$ws.Range("E26").Value = '  +1.97%  '
$ws.Range("E27").Value = '  -0.05%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.05'
$ws.Range("E28").Value = '  +4.63%  '
$ws.Range("E29").Value = '  -2.17%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '35.31'
$ws.Range("E30").Value = '  +7.07%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.135'
$ws.Range("E31").Value = '  +8.48%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '49.49'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '19.68'
$ws.Range("E33").Value = '  -3.29%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.39'
$ws.Range("E34").Value = '  +3.32%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0788'
$ws.Range("E35").Value = '  +2.80%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.00'
$ws.Range("E36").Value = '  +0.23%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.98'
$ws.Range("E37").Value = '  +6.12%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.64'
$ws.Range("E38").Value = '  +4.00%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.00'
$ws.Range("E39").Value = '  +4.37%  '
$ws.Range("E40").Value = '  +2.01%  '
$ws.Range("E41").Value = '  +1.59%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '121.68'
$ws.Range("E42").Value = '  -2.96%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '21.14'
$ws.Range("E43").Value = '  +2.69%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0298'
$ws.Range("E44").Value = '  +2.69%  '
$ws.Range("D45").Value = '1.968.20'
$ws.Range("E45").Value = '  +1.87%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.99'
$ws.Range("E46").Value = '  +2.51%  '
$ws.Range("E47").Value = '  -0.40%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.81'
$ws.Range("E48").Value = '  +1.23%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.13'
$ws.Range("E49").Value = '  -0.02%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '5.27'
$ws.Range("E50").Value = '  +11.35%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '79.69'
$ws.Range("E51").Value = '  +3.76%  '
